# Update CDA Logical model for ST.r2b
# Applies the metadata refresh (version bump, new date, Jurisdiction row)
# to the "Metadata" (Property/Value) worksheet of the SXCM-TS logical
# model workbook. The "Elements" worksheet is untouched content-wise;
# its shared-string indices shift automatically as a side effect of the
# shared string table changing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Version: 2.0.0-sd-202406-matchbox-patch -> 2.0.1-sd-202510-matchbox-patch
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# --- Date: 2024-06-19T17:47:42+02:00 -> 2025-10-29T22:15:57+01:00
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- Insert a new "Jurisdiction" row right after "Contact" (row 10),
#     pushing "Description" and everything below down by one row.
$ws.Rows.Item(11).Insert()
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# --- Append a new "Derivation" / "specialization" row at the end of the
#     table (row 20), matching the formatting of the row above it.
$ws.Range("A19:B19").Copy()
$ws.Range("A20:B20").PasteSpecial(-4122)
$ws.Range("A20").Value = "Derivation"
$ws.Range("B20").Value = "specialization"
